$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D cells are written as exact text (avoid Excel auto-converting
# numeric-looking strings like "240.85" into true numbers, which would lose
# trailing zeros / exact formatting).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '40.956.25'
$ws.Range("E2").Value = '  -6.62%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.183.35'
$ws.Range("E3").Value = '  -7.20%  '

# Row 4
$ws.Range("E4").Value = '  +0.12%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.85'
$ws.Range("E5").Value = '  +0.38%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.616'
$ws.Range("E6").Value = '  -7.59%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '69.16'
$ws.Range("E7").Value = '  -5.65%  '

# Row 8
$ws.Range("E8").Value = '  +0.34%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.535'
$ws.Range("E9").Value = '  -11.23%  '

# Row 10
$ws.Range("B10").Value = 'OKB'
$ws.Range("C10").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '57.68'
$ws.Range("E10").Value = '  -5.02%  '

# Row 11
$ws.Range("B11").Value = 'Avalanche'
$ws.Range("C11").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '36.21'
$ws.Range("E11").Value = '  +7.00%  '

# Row 12
$ws.Range("B12").Value = 'Dogecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0931'
$ws.Range("E12").Value = '  -8.46%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.104'
$ws.Range("E13").Value = '  -4.06%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.52'
$ws.Range("E14").Value = '  -9.85%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.520.74'
$ws.Range("E15").Value = '  -6.76%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.45'
$ws.Range("E16").Value = '  -10.62%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.827'
$ws.Range("E17").Value = '  -8.60%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.194.15'
$ws.Range("E18").Value = '  -6.74%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '40.987.53'
$ws.Range("E19").Value = '  -6.52%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0931'
$ws.Range("E20").Value = '  -9.61%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.95'
$ws.Range("E21").Value = '  -5.98%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.98'
$ws.Range("E22").Value = '  -8.33%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '230.01'
$ws.Range("E23").Value = '  -9.02%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.99'
$ws.Range("E24").Value = '  +7.08%  '

# Row 26
$ws.Range("E26").Value = '  -4.79%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.39'
$ws.Range("E27").Value = '  -3.67%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.18'
$ws.Range("E28").Value = '  -4.93%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.63'
$ws.Range("E29").Value = '  -7.70%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '167.12'
$ws.Range("E30").Value = '  -4.95%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '20.08'
$ws.Range("E31").Value = '  -9.58%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.117'
$ws.Range("E32").Value = '  -9.28%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.123'
$ws.Range("E33").Value = '  -8.09%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0695'
$ws.Range("E34").Value = '  -6.47%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.03'
$ws.Range("E35").Value = '  -5.73%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.52'
$ws.Range("E36").Value = '  -10.53%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.81'
$ws.Range("E37").Value = '  +0.79%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '23.29'
$ws.Range("E38").Value = '  +17.94%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.22'
$ws.Range("E39").Value = '  -7.40%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0265'
$ws.Range("E40").Value = '  -3.34%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.72'
$ws.Range("E41").Value = '  -11.75%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '64.38'
$ws.Range("E42").Value = '  -1.68%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.65'
$ws.Range("E43").Value = '  -4.49%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.67'
$ws.Range("E44").Value = '  -15.74%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.189'
$ws.Range("E45").Value = '  -5.24%  '

# Row 46
$ws.Range("E46").Value = '  +0.12%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0976'
$ws.Range("E47").Value = '  -8.34%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.44'
$ws.Range("E48").Value = '  +3.41%  '

# Row 49
$ws.Range("B49").Value = 'TrustWalletToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.15'
$ws.Range("E49").Value = '  -6.13%  '

# Row 50
$ws.Range("B50").Value = 'Celestia'
$ws.Range("C50").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.80'
$ws.Range("E50").Value = '  +1.69%  '

# Row 51
$ws.Range("E51").Value = '  -6.67%  '
